$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1756.9546
$ws.Range("I28").Value = 1419.375
$ws.Range("J28").Value = 2657.1667
$ws.Range("K28").Value = 1419.375
$ws.Range("L28").Value = 2657.1667
$ws.Range("M28").Value = -934.375
$ws.Range("N28").Value = -3627.1667
$ws.Range("H40").Value = 11846.88
$ws.Range("J40").Value = 13114.333
$ws.Range("L40").Value = 13114.333
$ws.Range("N40").Value = -13464.333
$ws.Range("H70").Value = 4530.067
$ws.Range("J70").Value = 5845.2
$ws.Range("L70").Value = 17535.6
$ws.Range("N70").Value = -18075.6
$ws.Range("H73").Value = 4530.067
$ws.Range("J73").Value = 5845.2
$ws.Range("L73").Value = 17535.6
$ws.Range("N73").Value = -19407.6
$ws.Range("H76").Value = 2800.5
$ws.Range("I76").Value = 2101.5
$ws.Range("J76").Value = 3499.5
$ws.Range("K76").Value = 2101.5
$ws.Range("L76").Value = 3499.5
$ws.Range("M76").Value = -1786.5
$ws.Range("N76").Value = -4129.5
$ws.Range("H79").Value = 2800.5
$ws.Range("I79").Value = 2101.5
$ws.Range("J79").Value = 3499.5
$ws.Range("K79").Value = 2101.5
$ws.Range("L79").Value = 3499.5
$ws.Range("M79").Value = -1009.5
$ws.Range("N79").Value = -5683.5
$ws.Range("H98").Value = 1031.48
$ws.Range("I98").Value = 1136.7727
$ws.Range("K98").Value = 1136.7727
$ws.Range("M98").Value = 361.2273
$ws.Range("H100").Value = 683.1429000000001
$ws.Range("I100").Value = 747.8333
$ws.Range("J100").Value = 295
$ws.Range("K100").Value = 747.8333
$ws.Range("L100").Value = 295
$ws.Range("M100").Value = -206.8333
$ws.Range("N100").Value = -1377
$ws.Range("H103").Value = 745.1539
$ws.Range("I103").Value = 730.4
$ws.Range("J103").Value = 794.3333
$ws.Range("K103").Value = 2191.2
$ws.Range("L103").Value = 2382.9999
$ws.Range("M103").Value = -1605.2
$ws.Range("N103").Value = -3554.9999
$ws.Range("H111").Value = 3008.9092
$ws.Range("I111").Value = 1757
$ws.Range("K111").Value = 5271
$ws.Range("M111").Value = -2204
$ws.Range("H112").Value = 1668.1111
$ws.Range("I112").Value = 902
$ws.Range("J112").Value = 1853.0344
$ws.Range("K112").Value = 2706
$ws.Range("L112").Value = 5559.1032
$ws.Range("M112").Value = -1598
$ws.Range("N112").Value = -7775.1032
$ws.Range("H122").Value = 1031.48
$ws.Range("I122").Value = 1136.7727
$ws.Range("K122").Value = 3410.3181
$ws.Range("M122").Value = -960.3181
$ws.Range("H129").Value = 2093.9
$ws.Range("I129").Value = 947.4167
$ws.Range("J129").Value = 2858.2222
$ws.Range("K129").Value = 2842.2501
$ws.Range("L129").Value = 8574.6666
$ws.Range("M129").Value = 2157.7499
$ws.Range("N129").Value = -18574.6666
$ws.Range("H132").Value = 12617.911
$ws.Range("I132").Value = 1728.2593
$ws.Range("J132").Value = 22756.55
$ws.Range("K132").Value = 5184.7779
$ws.Range("L132").Value = 68269.64999999999
$ws.Range("M132").Value = -2654.7779
$ws.Range("N132").Value = -73329.64999999999
$ws.Range("H137").Value = 11448403
$ws.Range("I137").Value = 715590.9
$ws.Range("J137").Value = 20839614
$ws.Range("K137").Value = 2146772.7
$ws.Range("L137").Value = 62518842
$ws.Range("M137").Value = -2144222.7
$ws.Range("N137").Value = -62523942
$ws.Range("H138").Value = 2162.28
$ws.Range("J138").Value = 2515.9453
$ws.Range("L138").Value = 7547.8359
$ws.Range("N138").Value = -17827.8359

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 700868.5
$ws.Range("I2").Value = 834062.5600000001
$ws.Range("K2").Value = 834062.5600000001
$ws.Range("M2").Value = -833949.5600000001
$ws.Range("H32").Value = 15191.889
$ws.Range("I32").Value = 15051.865
$ws.Range("J32").Value = 15853.818
$ws.Range("K32").Value = 15051.865
$ws.Range("L32").Value = 15853.818
$ws.Range("M32").Value = -14764.865
$ws.Range("N32").Value = -16427.818
$ws.Range("H45").Value = 2954.6365
$ws.Range("I45").Value = 1584.3334
$ws.Range("J45").Value = 4599
$ws.Range("K45").Value = 1584.3334
$ws.Range("L45").Value = 4599
$ws.Range("M45").Value = -1207.3334
$ws.Range("N45").Value = -5353
$ws.Range("H74").Value = 1305.3158
$ws.Range("I74").Value = 1223.7059
$ws.Range("K74").Value = 1223.7059
$ws.Range("M74").Value = -349.7058999999999
$ws.Range("H77").Value = 1305.3158
$ws.Range("I77").Value = 1223.7059
$ws.Range("K77").Value = 6118.5295
$ws.Range("M77").Value = -1750.5295
$ws.Range("H102").Value = 429406.3
$ws.Range("I102").Value = 548867.0600000001
$ws.Range("J102").Value = 2760.7144
$ws.Range("K102").Value = 548867.0600000001
$ws.Range("L102").Value = 2760.7144
$ws.Range("M102").Value = -547245.0600000001
$ws.Range("N102").Value = -6004.7144
$ws.Range("H116").Value = 700868.5
$ws.Range("I116").Value = 834062.5600000001
$ws.Range("K116").Value = 834062.5600000001
$ws.Range("M116").Value = -831768.5600000001
$ws.Range("H122").Value = 4398.8438
$ws.Range("I122").Value = 2494.9524
$ws.Range("K122").Value = 7484.8572
$ws.Range("M122").Value = -5034.8572
$ws.Range("H132").Value = 10383.614
$ws.Range("I132").Value = 12527.98
$ws.Range("J132").Value = 5022.7
$ws.Range("K132").Value = 37583.94
$ws.Range("L132").Value = 15068.1
$ws.Range("M132").Value = -35053.94
$ws.Range("N132").Value = -20128.1

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 700868.5
$ws.Range("I3").Value = 834062.5600000001
$ws.Range("K3").Value = 834062.5600000001
$ws.Range("M3").Value = -833948.5600000001
$ws.Range("H20").Value = 1657.9111
$ws.Range("I20").Value = 1209.1875
$ws.Range("J20").Value = 2762.4614
$ws.Range("K20").Value = 1209.1875
$ws.Range("L20").Value = 2762.4614
$ws.Range("M20").Value = -962.1875
$ws.Range("N20").Value = -3256.4614
$ws.Range("H58").Value = 37446.8
$ws.Range("I58").Value = 29420
$ws.Range("J58").Value = 39453.5
$ws.Range("K58").Value = 29420
$ws.Range("L58").Value = 39453.5
$ws.Range("M58").Value = -29126
$ws.Range("N58").Value = -40041.5
$ws.Range("H60").Value = 91311.8
$ws.Range("J60").Value = 89945
$ws.Range("L60").Value = 89945
$ws.Range("N60").Value = -91143
$ws.Range("H86").Value = 4533
$ws.Range("I86").Value = 3874.5
$ws.Range("K86").Value = 3874.5
$ws.Range("M86").Value = -2751.5
$ws.Range("H89").Value = 4533
$ws.Range("I89").Value = 3874.5
$ws.Range("K89").Value = 19372.5
$ws.Range("M89").Value = -13756.5
$ws.Range("H99").Value = 1097809.9
$ws.Range("I99").Value = 1489342
$ws.Range("K99").Value = 1489342
$ws.Range("M99").Value = -1487844
$ws.Range("H107").Value = 723.3333
$ws.Range("I107").Value = 723.5714
$ws.Range("J107").Value = 722.5
$ws.Range("K107").Value = 723.5714
$ws.Range("L107").Value = 722.5
$ws.Range("M107").Value = 1196.4286
$ws.Range("N107").Value = -4562.5
$ws.Range("H134").Value = 852.0571
$ws.Range("I134").Value = 832.5161000000001
$ws.Range("J134").Value = 1003.5
$ws.Range("K134").Value = 2497.5483
$ws.Range("L134").Value = 3010.5
$ws.Range("M134").Value = 37.45169999999962
$ws.Range("N134").Value = -8080.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6614.478
$ws.Range("I31").Value = 1774
$ws.Range("J31").Value = 8732.1875
$ws.Range("K31").Value = 1774
$ws.Range("L31").Value = 8732.1875
$ws.Range("M31").Value = -1479
$ws.Range("N31").Value = -9322.1875
$ws.Range("H34").Value = 6614.478
$ws.Range("I34").Value = 1774
$ws.Range("J34").Value = 8732.1875
$ws.Range("K34").Value = 1774
$ws.Range("L34").Value = 8732.1875
$ws.Range("M34").Value = -1572
$ws.Range("N34").Value = -9136.1875
$ws.Range("H36").Value = 20148
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H39").Value = 3500
$ws.Range("I39").Value = 3500
$ws.Range("K39").Value = 3500
$ws.Range("M39").Value = -3109
$ws.Range("H40").Value = 20148
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H49").Value = 3500
$ws.Range("I49").Value = 3500
$ws.Range("K49").Value = 3500
$ws.Range("M49").Value = -3318
$ws.Range("H58").Value = 373954.6
$ws.Range("I58").Value = 401731.16
$ws.Range("K58").Value = 401731.16
$ws.Range("M58").Value = -401528.16
$ws.Range("H86").Value = 10439
$ws.Range("I86").Value = 10313.5
$ws.Range("K86").Value = 10313.5
$ws.Range("M86").Value = -9190.5
$ws.Range("H89").Value = 10439
$ws.Range("I89").Value = 10313.5
$ws.Range("K89").Value = 51567.5
$ws.Range("M89").Value = -45951.5
$ws.Range("H107").Value = 1299301.9
$ws.Range("I107").Value = 1653428.5
$ws.Range("K107").Value = 1653428.5
$ws.Range("M107").Value = -1651508.5
$ws.Range("H132").Value = 9596.034
$ws.Range("I132").Value = 9596.034
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 28788.102
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -26258.102
$ws.Range("H134").Value = 1877.1482
$ws.Range("I134").Value = 1680.1154
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 5040.3462
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -2505.3462
$ws.Range("N134").Value = -26070
$ws.Range("H136").Value = 373954.6
$ws.Range("I136").Value = 401731.16
$ws.Range("K136").Value = 1205193.48
$ws.Range("M136").Value = -1202643.48
$ws.Range("H140").Value = 40000
$ws.Range("J140").Value = 40000
$ws.Range("L140").Value = 40000
$ws.Range("N140").Value = -50360
$ws.Range("H141").Value = 78693.12
$ws.Range("J141").Value = 79842.94
$ws.Range("L141").Value = 79842.94
$ws.Range("N141").Value = -90202.94

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 4999
$ws.Range("J17").Value = 4999
$ws.Range("L17").Value = 14997
$ws.Range("N17").Value = -15335
$ws.Range("H34").Value = 2053.1428
$ws.Range("J34").Value = 1467.6666
$ws.Range("L34").Value = 4402.9998
$ws.Range("N34").Value = -4570.9998
$ws.Range("H39").Value = 5657.357
$ws.Range("J39").Value = 5712.154
$ws.Range("L39").Value = 17136.462
$ws.Range("N39").Value = -17724.462
$ws.Range("H51").Value = 2640.6
$ws.Range("J51").Value = 2799.75
$ws.Range("L51").Value = 8399.25
$ws.Range("N51").Value = -9319.25
$ws.Range("H55").Value = 2370.9412
$ws.Range("J55").Value = 3273.4546
$ws.Range("L55").Value = 9820.363799999999
$ws.Range("N55").Value = -10174.3638
$ws.Range("H57").Value = 23333.334
$ws.Range("I57").Value = 20000
$ws.Range("K57").Value = 60000
$ws.Range("M57").Value = -59441
$ws.Range("H103").Value = 800
$ws.Range("I103").Value = 800
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 2400
$ws.Range("L103").ClearContents()
$ws.Range("N103").Value = 0
$ws.Range("M103").Value = -1521
$ws.Range("H137").Value = 77286900
$ws.Range("I137").Value = 93754490
$ws.Range("K137").Value = 281263470
$ws.Range("M137").Value = -281258370

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4137450
$ws.Range("I70").Value = 6498064.5
$ws.Range("J70").Value = 6374.5
$ws.Range("K70").Value = 6498064.5
$ws.Range("L70").Value = 6374.5
$ws.Range("M70").Value = -6497794.5
$ws.Range("N70").Value = -6914.5
$ws.Range("H73").Value = 4137450
$ws.Range("I73").Value = 6498064.5
$ws.Range("J73").Value = 6374.5
$ws.Range("K73").Value = 6498064.5
$ws.Range("L73").Value = 6374.5
$ws.Range("M73").Value = -6497128.5
$ws.Range("N73").Value = -8246.5
$ws.Range("H80").Value = 1210828.6
$ws.Range("I80").Value = 2378414.8
$ws.Range("J80").Value = 43242.57
$ws.Range("K80").Value = 2378414.8
$ws.Range("L80").Value = 43242.57
$ws.Range("M80").Value = -2377416.8
$ws.Range("N80").Value = -45238.57
$ws.Range("H83").Value = 1210828.6
$ws.Range("I83").Value = 2378414.8
$ws.Range("J83").Value = 43242.57
$ws.Range("K83").Value = 11892074
$ws.Range("L83").Value = 216212.85
$ws.Range("M83").Value = -11887082
$ws.Range("N83").Value = -226196.85
$ws.Range("H94").Value = 62500
$ws.Range("J94").Value = 60000
$ws.Range("L94").Value = 60000
$ws.Range("N94").Value = -61352
$ws.Range("H102").Value = 10708.632
$ws.Range("I102").Value = 11420.538
$ws.Range("K102").Value = 11420.538
$ws.Range("M102").Value = -9798.538
$ws.Range("H122").Value = 922096.8
$ws.Range("I122").Value = 1836201.1
$ws.Range("K122").Value = 5508603.300000001
$ws.Range("M122").Value = -5506153.300000001
$ws.Range("H132").Value = 576750.9399999999
$ws.Range("I132").Value = 226208.44
$ws.Range("J132").Value = 839657.8
$ws.Range("K132").Value = 678625.3200000001
$ws.Range("L132").Value = 2518973.4
$ws.Range("M132").Value = -676095.3200000001
$ws.Range("N132").Value = -2524033.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3022.08
$ws.Range("I16").Value = 1652.85
$ws.Range("K16").Value = 1652.85
$ws.Range("M16").Value = -1482.85
$ws.Range("H22").Value = 1109.1818
$ws.Range("I22").Value = 950.1667
$ws.Range("J22").Value = 1300
$ws.Range("K22").Value = 950.1667
$ws.Range("L22").Value = 1300
$ws.Range("M22").Value = -655.1667
$ws.Range("N22").Value = -1890
$ws.Range("H27").Value = 1109.1818
$ws.Range("I27").Value = 950.1667
$ws.Range("J27").Value = 1300
$ws.Range("K27").Value = 950.1667
$ws.Range("L27").Value = 1300
$ws.Range("M27").Value = -843.1667
$ws.Range("N27").Value = -1514
$ws.Range("H40").Value = 1575
$ws.Range("I40").Value = 1575
$ws.Range("K40").Value = 1575
$ws.Range("M40").Value = -1439
$ws.Range("H46").Value = 5991.48
$ws.Range("I46").Value = 4600
$ws.Range("J46").Value = 6049.4585
$ws.Range("K46").Value = 4600
$ws.Range("L46").Value = 6049.4585
$ws.Range("M46").Value = -4412
$ws.Range("N46").Value = -6425.4585
$ws.Range("H55").Value = 366
$ws.Range("I55").Value = 159
$ws.Range("J55").Value = 407.4
$ws.Range("K55").Value = 159
$ws.Range("L55").Value = 407.4
$ws.Range("M55").Value = 14
$ws.Range("N55").Value = -753.4
$ws.Range("H61").Value = 3201.8572
$ws.Range("I61").Value = 603.5
$ws.Range("K61").Value = 603.5
$ws.Range("M61").Value = -401.5
$ws.Range("H68").Value = 2846282.5
$ws.Range("I68").Value = 22727272
$ws.Range("J68").Value = 6141.143
$ws.Range("K68").Value = 22727272
$ws.Range("L68").Value = 6141.143
$ws.Range("M68").Value = -22726523
$ws.Range("N68").Value = -7639.143
$ws.Range("H71").Value = 2846282.5
$ws.Range("I71").Value = 22727272
$ws.Range("J71").Value = 6141.143
$ws.Range("K71").Value = 113636360
$ws.Range("L71").Value = 30705.715
$ws.Range("M71").Value = -113632616
$ws.Range("N71").Value = -38193.715
$ws.Range("H113").Value = 3201.8572
$ws.Range("I113").Value = 603.5
$ws.Range("K113").Value = 603.5
$ws.Range("M113").Value = 1566.5
$ws.Range("H127").Value = 250095070
$ws.Range("J127").Value = 126770.664
$ws.Range("L127").Value = 126770.664
$ws.Range("N127").Value = -136690.664
$ws.Range("H132").Value = 4826.4346
$ws.Range("I132").Value = 4639.3335
$ws.Range("K132").Value = 13918.0005
$ws.Range("M132").Value = -11388.0005
$ws.Range("H134").Value = 85428.5
$ws.Range("J134").Value = 85428.5
$ws.Range("L134").Value = 85428.5
$ws.Range("N134").Value = -95568.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H107").Value = 4539.5835
$ws.Range("I107").Value = 5111.9414
$ws.Range("K107").Value = 15335.8242
$ws.Range("M107").Value = -13415.8242
$ws.Range("H122").Value = 5718.1177
$ws.Range("I122").Value = 5575.5625
$ws.Range("K122").Value = 16726.6875
$ws.Range("M122").Value = -14276.6875
$ws.Range("H136").Value = 8247.794
$ws.Range("I136").Value = 2970.125
$ws.Range("K136").Value = 8910.375
$ws.Range("M136").Value = -6360.375
